$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 808
$ws.Range("I2").Value = 2083
$ws.Range("J2").Value = 8571
$ws.Range("K2").Value = 52
$ws.Range("L2").Value = 2323
$ws.Range("M2").Value = 131
$ws.Range("N2").Value = 1526
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 35
$ws.Range("Q2").Value = 16
$ws.Range("R2").Value = 101
$ws.Range("S2").Value = 938
$ws.Range("T2").Value = 1577
$ws.Range("U2").Value = 93
$ws.Range("V2").Value = 13439
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 13348
$ws.Range("Y2").Value = 18
$ws.Range("Z2").Value = 226
$ws.Range("AA2").Value = 87
